# Edit script for Työaikaraportti_ArttuMutka.xlsx
# Implements:
#  - a new timesheet row (45347 / 6h / "about page" entry) inserted before the totals row
#  - the existing last entry's description gets ". Aloitin systeminformation." appended
#  - the totals row moves down one row and its SUM formula / result updates accordingly
#  - workbook/sheet view (window size, scroll position, selection) updated to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Shift the totals row (old row 25) down to row 26, then turn old row 24
#    (last data row) into the template for the new row 25 by copying its
#    formatting down. Using targeted Range->Range copies (rather than
#    Rows.Insert) keeps the existing style table intact instead of growing it.
# ---------------------------------------------------------------------------
$ws.Range("B25:D25").Copy($ws.Range("B26:D26"))
$ws.Range("B24:D24").Copy($ws.Range("B25:D25"))

# Row heights: new data row (25) matches the other multi-line entries (37.5),
# the relocated totals row (26) keeps its original height (18.75).
$ws.Rows("25:25").RowHeight = 37.5
$ws.Rows("26:26").RowHeight = 18.75

# ---------------------------------------------------------------------------
# 2) Fill in the brand-new row 25 with the new timesheet entry. D25 still
#    shares its text with D24 at this point (both came from the row-24
#    copy), so setting it first makes the new "Paransin..." text reuse that
#    original shared-string slot while D24 keeps pointing at the old text.
# ---------------------------------------------------------------------------
$ws.Range("B25").Value = 45347
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = "Paransin about sivua. Myöskin koitin ideoida ja parantaa systeminformation työkalua"

# ---------------------------------------------------------------------------
# 3) Update the text of the (old) last entry, row 24, appending the new
#    sentence about starting on systeminformation. D24 is now the sole
#    reference to the old text, so this edits that shared-string entry
#    in place.
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = "Aloin suunitelemaana lisää ominaisuuksia. Implementoin uusia usercontrol ja niille omat navigaatio osiot. Lajitelin tiedostoja paremmin. Tein random number generaatorin. Tein Yksinkertaisen checksum työkalun. Aloitin systeminformation."

# ---------------------------------------------------------------------------
# 4) Fix up the totals row, now on row 26: label stays the same, formula sums
#    the new range of data rows (C6:C25) and recalculates the total.
# ---------------------------------------------------------------------------
$ws.Range("C26").Formula = "=SUM(C6:C25)"

# ---------------------------------------------------------------------------
# 5) View/window bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$w = $excel.ActiveWindow
$w.Left = 2220
$w.Top = 630
$w.Width = 21600
$w.Height = 11295
$w.ScrollRow = 22
$w.ScrollColumn = 1

$ws.Range("I29").Select()
